$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell G1 = "Problem" (shared string index 13), styled like
#     the rest of the header row (s="2") ---
$ws.Range("G1").Value = "Problem"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# --- Updated progress numbers ---
$ws.Range("E2").Value = 0.2
$ws.Range("E5").Value = 0.85

# --- New text cells, set in an order that reproduces the shared-string
#     table order seen in the target file: Problem, Last test,
#     "the whole subject", Ensemble learning ---
$ws.Range("G5").Value = "Last test"
$ws.Range("G2").Value = "the whole subject"

# D5 must hold the existing "2025.01.12" text (not be auto-converted to a
# date), and keep the default (unstyled) cell format just like D3/D4, so
# copy the value+format from D3 which already contains exactly that text.
$ws.Range("D3").Copy()
$ws.Range("D5").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("F5").Value = "Ensemble learning"

$excel.CutCopyMode = 0

# --- Selection shown in the diff ---
$ws.Range("F5").Select()
